$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "html" file-format entry, mirroring the existing txt/xml/json
# blocks (one row per language: ara, eng, fra) at the bottom of the table.

# ara row (12)
$ws.Range("A12").Value = "html"
$ws.Range("B12").Value = "ملف html"
$ws.Range("C12").Value = "ara"
$ws.Range("D12").Value = $true
$ws.Range("D12").HorizontalAlignment = -4131
$ws.Range("E12").Value = "superadmin"
$ws.Range("F12").Value = "now()"

# eng row (11)
$ws.Range("A11").Value = "html"
$ws.Range("B11").Value = "html file"
$ws.Range("C11").Value = "eng"
$ws.Range("D11").Value = $true
$ws.Range("D11").HorizontalAlignment = -4131
$ws.Range("E11").Value = "superadmin"
$ws.Range("F11").Value = "now()"

# fra row (13)
$ws.Range("A13").Value = "html"
$ws.Range("B13").Value = "Fichier html"
$ws.Range("C13").Value = "fra"
$ws.Range("D13").Value = $true
$ws.Range("D13").HorizontalAlignment = -4131
$ws.Range("E13").Value = "superadmin"
$ws.Range("F13").Value = "now()"

# Mirror the author's final selection state (column G selected).
$ws.Range("G1:XFD1048576").Select() | Out-Null
